$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 15:21"

# Refresh COVID data values (table stays sorted desc by column B; only cell values change)
# Row 4
$ws.Range("B4").Value = 2235678
$ws.Range("C4").Value = 1207
$ws.Range("E4").Value = 1196927
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 119955
# Row 14
$ws.Range("B14").Value = 189504
$ws.Range("E14").Value = 6977
# Row 19
$ws.Range("B19").Value = 145991
$ws.Range("C19").Value = 4757
$ws.Range("D19").Value = 93915
$ws.Range("E19").Value = 50937
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 1139
# Row 22
$ws.Range("B22").Value = 84441
$ws.Range("C22").Value = 1267
$ws.Range("D22").Value = 63642
$ws.Range("E22").Value = 20713
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 86
# Row 23
$ws.Range("B23").Value = 83293
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 78394
$ws.Range("E23").Value = 265
$ws.Range("H23").Value = 4634
# Row 28
$ws.Range("B28").Value = 56043
$ws.Range("C28").Value = 1481
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 5053
# Row 29
$ws.Range("B29").Value = 49319
$ws.Range("C29").Value = 115
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 6078
# Row 30
$ws.Range("B30").Value = 49219
$ws.Range("D30").Value = 13141
$ws.Range("E30").Value = 34228
$ws.Range("H30").Value = 1850
# Row 35
$ws.Range("B35").Value = 38089
$ws.Range("C35").Value = 417
$ws.Range("D35").Value = 24010
$ws.Range("E35").Value = 12555
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1524
# Row 36
$ws.Range("B36").Value = 38074
$ws.Range("C36").Value = 541
$ws.Range("D36").Value = 29512
$ws.Range("E36").Value = 8254
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 308
# Row 37
$ws.Range("D37").Value = 10721
$ws.Range("E37").Value = 23902
$ws.Range("G37").Value = 16
$ws.Range("H37").Value = 929
# Row 56
$ws.Range("D56").Value = 10065
$ws.Range("E56").Value = 5712
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 100
# Row 76
$ws.Range("B76").Value = 5730
$ws.Range("C76").Value = 48
$ws.Range("D76").Value = 4166
$ws.Range("E76").Value = 1545
# Row 77
$ws.Range("B77").Value = 5475
$ws.Range("C77").Value = 106
$ws.Range("D77").Value = 3716
$ws.Range("E77").Value = 1683
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 76
# Row 82
$ws.Range("B82").Value = 4664
$ws.Range("C82").Value = 182
$ws.Range("D82").Value = 1836
$ws.Range("E82").Value = 2612
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 216
# Row 83
$ws.Range("B83").Value = 4545
$ws.Range("D83").Value = 3411
$ws.Range("E83").Value = 1091
$ws.Range("H83").Value = 43
# Row 89
$ws.Range("B89").Value = 3954
$ws.Range("C89").Value = 195
$ws.Range("D89").Value = 934
$ws.Range("E89").Value = 2955
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 65
# Row 93
$ws.Range("B93").Value = 3174
$ws.Range("C93").Value = 33
$ws.Range("D93").Value = 2219
$ws.Range("E93").Value = 787
# Row 100
$ws.Range("B100").Value = 2269
$ws.Range("C100").Value = 11
$ws.Range("D100").Value = 2142
$ws.Range("E100").Value = 20
# Row 108
$ws.Range("B108").Value = 1816
$ws.Range("C108").Value = 1
$ws.Range("E108").Value = 9
# Row 117
$ws.Range("B117").Value = 1495
$ws.Range("C117").Value = 6
$ws.Range("D117").Value = 944
$ws.Range("E117").Value = 519
$ws.Range("H117").Value = 32
# Row 118
$ws.Range("B118").Value = 1492
$ws.Range("D118").Value = 153
$ws.Range("E118").Value = 1324
$ws.Range("H118").Value = 15
# Row 145
$ws.Range("B145").Value = 597
$ws.Range("C145").Value = 25
$ws.Range("D145").Value = 238
$ws.Range("E145").Value = 348
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 11
# Row 146
$ws.Range("B146").Value = 579
$ws.Range("C146").Value = 24
$ws.Range("D146").Value = 415
$ws.Range("E146").Value = 161
$ws.Range("H146").Value = 3
# Row 147
$ws.Range("D147").Value = 73
$ws.Range("E147").Value = 493
$ws.Range("H147").Value = 6
# Row 208
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0
# Row 209
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
